# Bulk orders landing page: quantities were populated for a few line items
# and the mobile-view header row (row 18) is left selected/in-focus, matching
# the last edited cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 2
$ws.Range("A12").Value = 2
$ws.Range("A18").Value = 3

$ws.Range("A18").Select()
